# [ANV] updated the tapes for 1ppb
#
# For each isotope-detail row in the "Target Fractions" sheet, populate
# column B with the parent nuclide's abundance/branching value (copied
# down from the nearest preceding row that already carries a B value),
# and rewrite the F/G/H formulas so that they are additionally scaled
# by that column-B value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Target Fractions")

# Map of row number -> value to place into column B of that row.
# (values written in plain decimal form, no scientific notation)
$bValues = @{
    4  = 0.00070227284042469250
    6  = 0.63266447300513628527
    7  = 0.63266447300513628527
    9  = 0.02063825428636016993
    11 = 0.01023791290051672027
    12 = 0.01023791290051672027
    13 = 0.01023791290051672027
    15 = 0.04720578993057022998
    17 = 0.20948218342998464747
    18 = 0.20948218342998464747
    19 = 0.20948218342998464747
    21 = 0.00949248329975991602
    22 = 0.00949248329975991602
    23 = 0.00949248329975991602
    25 = 0.05019580024183650790
    26 = 0.05019580024183650790
    27 = 0.05019580024183650790
    28 = 0.05019580024183650790
    29 = 0.05019580024183650790
    31 = 0.00142021331426823430
    33 = 0.00959123362848334285
    34 = 0.00959123362848334285
    35 = 0.00959123362848334285
    36 = 0.00959123362848334285
}

$rows = @(4,6,7,9,11,12,13,15,17,18,19,21,22,23,25,26,27,28,29,31,33,34,35,36)

foreach ($r in $rows) {
    $b = $bValues[$r]
    $ws.Range("B$r").Value = $b
    $ws.Range("F$r").Formula = "=B$r*(D$r+E$r)/200"
    $ws.Range("G$r").Formula = "=B$r*D$r/100"
    $ws.Range("H$r").Formula = "=B$r*E$r/100"
}

$wb.Save()
